$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G header
$ws.Range("G1").Value = "Email-id"

# G2: email address with hyperlink formatting (mailto link)
$ws.Hyperlinks.Add($ws.Range("G2"), "mailto:alice@gmail.com", "", "", "alice@gmail.com")

# Row 4: a new set of (invalid) profile data entered for verification
$ws.Range("A4").Value = 9176226906
$ws.Range("B4").Value = "12Alice"
$ws.Range("C4").Value = '$K'
$ws.Range("D4").Value = "99-99-9999"
$ws.Range("G4").Value = "alice.com"

[void]$ws.Range("I5").Select()
